$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 183234.8
$ws.Range("I38").Value = 183234.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 549704.3999999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -549332.3999999999
$ws.Range("N38").ClearContents()

$ws.Range("H58").Value = 3991.111
$ws.Range("I58").Value = 390.25
$ws.Range("J58").Value = 6871.8
$ws.Range("K58").Value = 1170.75
$ws.Range("L58").Value = 20615.4
$ws.Range("M58").Value = -1020.75
$ws.Range("N58").Value = -20915.4

$ws.Range("H86").Value = 4788907
$ws.Range("J86").Value = 5852608.5
$ws.Range("L86").Value = 5852608.5
$ws.Range("N86").Value = -5854854.5

$ws.Range("H89").Value = 4788907
$ws.Range("J89").Value = 5852608.5
$ws.Range("L89").Value = 29263042.5
$ws.Range("N89").Value = -29274274.5

$ws.Range("H132").Value = 1743.4186
$ws.Range("I132").Value = 1594.4524
$ws.Range("K132").Value = 4783.357199999999
$ws.Range("M132").Value = -2253.357199999999

$ws.Range("H138").Value = 3944.0547
$ws.Range("I138").Value = 2304.7778
$ws.Range("J138").Value = 4174.578
$ws.Range("K138").Value = 6914.3334
$ws.Range("L138").Value = 12523.734
$ws.Range("M138").Value = -1774.3334
$ws.Range("N138").Value = -22803.734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3613.647
$ws.Range("I32").Value = 3280.506
$ws.Range("K32").Value = 3280.506
$ws.Range("M32").Value = -2993.506

$ws.Range("H45").Value = 3433.3635
$ws.Range("I45").Value = 2252.7144
$ws.Range("K45").Value = 2252.7144
$ws.Range("M45").Value = -1875.7144

$ws.Range("H113").Value = 99999
$ws.Range("J113").Value = 99999
$ws.Range("L113").Value = 99999
$ws.Range("N113").Value = -108677

$ws.Range("H122").Value = 5186.92
$ws.Range("I122").Value = 5859
$ws.Range("K122").Value = 17577
$ws.Range("M122").Value = -15127

$ws.Range("H124").Value = 79993
$ws.Range("J124").Value = 79993
$ws.Range("L124").Value = 79993
$ws.Range("N124").Value = -89813

$ws.Range("H125").Value = 71464
$ws.Range("J125").Value = 71464
$ws.Range("L125").Value = 71464
$ws.Range("N125").Value = -81304

$ws.Range("H132").Value = 2562.9148
$ws.Range("I132").Value = 2700.0698
$ws.Range("K132").Value = 8100.209400000001
$ws.Range("M132").Value = -5570.209400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 100237.2
$ws.Range("J74").Value = 100237.2
$ws.Range("L74").Value = 100237.2
$ws.Range("N74").Value = -102109.2

$ws.Range("H77").Value = 100237.2
$ws.Range("J77").Value = 100237.2
$ws.Range("L77").Value = 300711.6
$ws.Range("N77").Value = -310071.6

$ws.Range("H86").Value = 1419182.1
$ws.Range("I86").Value = 1890909.6
$ws.Range("K86").Value = 1890909.6
$ws.Range("M86").Value = -1889786.6

$ws.Range("H89").Value = 1419182.1
$ws.Range("I89").Value = 1890909.6
$ws.Range("K89").Value = 9454548
$ws.Range("M89").Value = -9448932

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 461.92856
$ws.Range("I107").Value = 393.0909
$ws.Range("J107").Value = 714.3333
$ws.Range("K107").Value = 393.0909
$ws.Range("L107").Value = 714.3333
$ws.Range("M107").Value = 1526.9091
$ws.Range("N107").Value = -4554.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 315.625
$ws.Range("I2").Value = 59.333332
$ws.Range("K2").Value = 355.999992
$ws.Range("M2").Value = -242.999992

$ws.Range("J4").Value = 350032.34
$ws.Range("L4").Value = 1050097.02
$ws.Range("N4").Value = -1050321.02

$ws.Range("H7").Value = 573
$ws.Range("I7").Value = 193.33333
$ws.Range("J7").Value = 762.8333
$ws.Range("K7").Value = 579.99999
$ws.Range("L7").Value = 2288.4999
$ws.Range("M7").Value = -467.99999
$ws.Range("N7").Value = -2512.4999

$ws.Range("H46").Value = 200
$ws.Range("I46").Value = 200
$ws.Range("K46").Value = 600
$ws.Range("M46").Value = -509

$ws.Range("H56").Value = 6583
$ws.Range("I56").Value = 6583
$ws.Range("K56").Value = 6583
$ws.Range("M56").Value = -6053

$ws.Range("H68").Value = 4000655.2
$ws.Range("I68").Value = 5000401
$ws.Range("J68").Value = 3334158
$ws.Range("K68").Value = 15001203
$ws.Range("L68").Value = 10002474
$ws.Range("M68").Value = -15000392
$ws.Range("N68").Value = -10004096

$ws.Range("H71").Value = 4000655.2
$ws.Range("I71").Value = 5000401
$ws.Range("J71").Value = 3334158
$ws.Range("K71").Value = 45003609
$ws.Range("L71").Value = 30007422
$ws.Range("M71").Value = -44999553
$ws.Range("N71").Value = -30015534

$ws.Range("H113").Value = 2180414
$ws.Range("J113").Value = 2236.6365
$ws.Range("L113").Value = 6709.9095
$ws.Range("N113").Value = -11049.9095

$ws.Range("H131").Value = 12424569
$ws.Range("I131").Value = 111445110
$ws.Range("J131").Value = 47001.625
$ws.Range("K131").Value = 334335330
$ws.Range("L131").Value = 141004.875
$ws.Range("M131").Value = -334330290
$ws.Range("N131").Value = -151084.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 6808.4
$ws.Range("I99").Value = 6808.4
$ws.Range("K99").Value = 6808.4
$ws.Range("M99").Value = -4562.4

$ws.Range("H122").Value = 4701.5557
$ws.Range("I122").Value = 2503.5
$ws.Range("K122").Value = 7510.5
$ws.Range("M122").Value = -5060.5

$ws.Range("H132").Value = 54126.19
$ws.Range("I132").Value = 6309.4
$ws.Range("K132").Value = 18928.2
$ws.Range("M132").Value = -16398.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7500
$ws.Range("I7").Value = 8588.888999999999
$ws.Range("K7").Value = 8588.888999999999
$ws.Range("M7").Value = -8476.888999999999

$ws.Range("H50").Value = 116080
$ws.Range("J50").Value = 116080
$ws.Range("L50").Value = 116080
$ws.Range("N50").Value = -117354

$ws.Range("H68").Value = 3617.6667
$ws.Range("J68").Value = 3926.5
$ws.Range("L68").Value = 3926.5
$ws.Range("N68").Value = -5424.5

$ws.Range("H71").Value = 3617.6667
$ws.Range("J71").Value = 3926.5
$ws.Range("L71").Value = 19632.5
$ws.Range("N71").Value = -27120.5

$ws.Range("H126").Value = 7500
$ws.Range("I126").Value = 8588.888999999999
$ws.Range("K126").Value = 25766.667
$ws.Range("M126").Value = -23296.667

$ws.Range("H132").Value = 7204.4375
$ws.Range("I132").Value = 6662.5713
$ws.Range("J132").Value = 10997.5
$ws.Range("K132").Value = 19987.7139
$ws.Range("L132").Value = 32992.5
$ws.Range("M132").Value = -17457.7139
$ws.Range("N132").Value = -38052.5

$ws.Range("H136").Value = 1005360
$ws.Range("I136").Value = 1433099.2
$ws.Range("J136").Value = 7301.6665
$ws.Range("K136").Value = 4299297.6
$ws.Range("L136").Value = 21904.9995
$ws.Range("M136").Value = -4296747.6
$ws.Range("N136").Value = -27004.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8552
$ws.Range("I74").Value = 5997
$ws.Range("K74").Value = 5997
$ws.Range("M74").Value = -5061

$ws.Range("H77").Value = 8552
$ws.Range("I77").Value = 5997
$ws.Range("K77").Value = 17991
$ws.Range("M77").Value = -13311

$ws.Range("H132").Value = 28187.45
$ws.Range("I132").Value = 2258.6
$ws.Range("K132").Value = 6775.799999999999
$ws.Range("M132").Value = -4245.799999999999

$ws.Range("H136").Value = 9335327
$ws.Range("I136").Value = 12718464
$ws.Range("K136").Value = 38155392
$ws.Range("M136").Value = -38152842
